# Update Brynhildr_Profits leve-crafting profit figures (scheduled market-price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7617.2666
$ws.Range("J132").Value = 750
$ws.Range("L132").Value = 2250
$ws.Range("N132").Value = -7310
$ws.Range("H137").Value = 12727.091
$ws.Range("I137").Value = 3999.125
$ws.Range("K137").Value = 11997.375
$ws.Range("M137").Value = -9447.375
$ws.Range("H138").Value = 4335.971
$ws.Range("I138").Value = 6522.9
$ws.Range("J138").Value = 3461.2
$ws.Range("K138").Value = 19568.7
$ws.Range("L138").Value = 10383.6
$ws.Range("M138").Value = -14428.7
$ws.Range("N138").Value = -20663.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2138.75
$ws.Range("I122").Value = 2302
$ws.Range("K122").Value = 6906
$ws.Range("M122").Value = -4456
$ws.Range("H132").Value = 4954.2065
$ws.Range("I132").Value = 3172.439
$ws.Range("K132").Value = 9517.316999999999
$ws.Range("M132").Value = -6987.316999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 62772
$ws.Range("I20").Value = 95365.09
$ws.Range("J20").Value = 3018
$ws.Range("K20").Value = 95365.09
$ws.Range("L20").Value = 3018
$ws.Range("M20").Value = -95118.09
$ws.Range("N20").Value = -3512
$ws.Range("H99").Value = 7170.5557
$ws.Range("I99").Value = 9907
$ws.Range("J99").Value = 1697.6666
$ws.Range("K99").Value = 9907
$ws.Range("L99").Value = 1697.6666
$ws.Range("M99").Value = -8409
$ws.Range("N99").Value = -4693.6666
$ws.Range("H100").Value = 24250
$ws.Range("J100").Value = 24250
$ws.Range("L100").Value = 24250
$ws.Range("N100").Value = -26414
$ws.Range("H107").Value = 1540.625
$ws.Range("I107").Value = 1387.5
$ws.Range("K107").Value = 1387.5
$ws.Range("M107").Value = 532.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2425.5881
$ws.Range("I31").Value = 3523.625
$ws.Range("J31").Value = 1449.5555
$ws.Range("K31").Value = 3523.625
$ws.Range("L31").Value = 1449.5555
$ws.Range("M31").Value = -3228.625
$ws.Range("N31").Value = -2039.5555
$ws.Range("H34").Value = 2425.5881
$ws.Range("I34").Value = 3523.625
$ws.Range("J34").Value = 1449.5555
$ws.Range("K34").Value = 3523.625
$ws.Range("L34").Value = 1449.5555
$ws.Range("M34").Value = -3321.625
$ws.Range("N34").Value = -1853.5555
$ws.Range("H99").Value = 12997.579
$ws.Range("I99").Value = 20451.727
$ws.Range("J99").Value = 2748.125
$ws.Range("K99").Value = 20451.727
$ws.Range("L99").Value = 2748.125
$ws.Range("M99").Value = -18953.727
$ws.Range("N99").Value = -5744.125
$ws.Range("H106").Value = 16083.25
$ws.Range("I106").Value = 15000
$ws.Range("J106").Value = 18249.75
$ws.Range("K106").Value = 15000
$ws.Range("L106").Value = 18249.75
$ws.Range("M106").Value = -13738
$ws.Range("N106").Value = -20773.75
$ws.Range("H126").Value = 12997.579
$ws.Range("I126").Value = 20451.727
$ws.Range("J126").Value = 2748.125
$ws.Range("K126").Value = 61355.181
$ws.Range("L126").Value = 8244.375
$ws.Range("M126").Value = -58885.181
$ws.Range("N126").Value = -13184.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 749
$ws.Range("J2").Value = 911.5
$ws.Range("L2").Value = 5469
$ws.Range("N2").Value = -5695
$ws.Range("H5").Value = 1244.9445
$ws.Range("J5").Value = 1764
$ws.Range("L5").Value = 5292
$ws.Range("N5").Value = -5516
$ws.Range("H11").Value = 142857380
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1331
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 1500
$ws.Range("M27").Value = -1398
$ws.Range("H80").Value = 15770.667
$ws.Range("J80").Value = 15770.667
$ws.Range("L80").Value = 47312.001
$ws.Range("N80").Value = -49184.001
$ws.Range("H83").Value = 15770.667
$ws.Range("J83").Value = 15770.667
$ws.Range("L83").Value = 141936.003
$ws.Range("N83").Value = -151296.003
$ws.Range("H135").Value = 1244.9445
$ws.Range("J135").Value = 1764
$ws.Range("L135").Value = 15876
$ws.Range("N135").Value = -20946

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3838.5715
$ws.Range("I122").Value = 3976.3635
$ws.Range("J122").Value = 3333.3333
$ws.Range("K122").Value = 11929.0905
$ws.Range("L122").Value = 9999.999899999999
$ws.Range("M122").Value = -9479.0905
$ws.Range("N122").Value = -14899.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5578.273
$ws.Range("I7").Value = 5040.8887
$ws.Range("K7").Value = 5040.8887
$ws.Range("M7").Value = -4928.8887
$ws.Range("H61").Value = 8103.6665
$ws.Range("I61").Value = 8108.727
$ws.Range("J61").Value = 8075.8335
$ws.Range("K61").Value = 8108.727
$ws.Range("L61").Value = 8075.8335
$ws.Range("M61").Value = -7906.727
$ws.Range("N61").Value = -8479.833500000001
$ws.Range("H93").Value = 1452.2142
$ws.Range("J93").Value = 2997.25
$ws.Range("L93").Value = 2997.25
$ws.Range("N93").Value = -5493.25
$ws.Range("H101").Value = 14125
$ws.Range("J101").Value = 14125
$ws.Range("L101").Value = 14125
$ws.Range("N101").Value = -20615
$ws.Range("H113").Value = 8103.6665
$ws.Range("I113").Value = 8108.727
$ws.Range("J113").Value = 8075.8335
$ws.Range("K113").Value = 8108.727
$ws.Range("L113").Value = 8075.8335
$ws.Range("M113").Value = -5938.727
$ws.Range("N113").Value = -12415.8335
$ws.Range("H126").Value = 5578.273
$ws.Range("I126").Value = 5040.8887
$ws.Range("K126").Value = 15122.6661
$ws.Range("M126").Value = -12652.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1492.5778
$ws.Range("I107").Value = 905.1
$ws.Range("K107").Value = 2715.3
$ws.Range("M107").Value = -795.3000000000002
$ws.Range("H122").Value = 33172.973
$ws.Range("I122").Value = 2856.3333
$ws.Range("J122").Value = 163101.42
$ws.Range("K122").Value = 8568.999899999999
$ws.Range("L122").Value = 489304.26
$ws.Range("M122").Value = -6118.999899999999
$ws.Range("N122").Value = -494204.26
$ws.Range("H126").Value = 1561.44
$ws.Range("I126").Value = 1380.35
$ws.Range("K126").Value = 4141.049999999999
$ws.Range("M126").Value = -1671.049999999999
$ws.Range("H136").Value = 1327.4717
$ws.Range("I136").Value = 1141.762
$ws.Range("K136").Value = 3425.286
$ws.Range("M136").Value = -875.2860000000001
